# Add 7 new rows of data (rows 3-9) to WorkSheet 1, mirroring the
# existing row 2 pattern (same Method/word-count style data), and widen
# column A slightly to fit the new content.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Each new row keeps the same Method/word-count style values as row 2;
# only the Date (col A) and ElapsedMs (col C) columns change per row.
$rows = @(
    @(42601.767372685186, 3233),
    @(42601.769293981481, 3141),
    @(42601.770046296297, 3085),
    @(42601.771041666667, 3000),
    @(42601.771458333336, 3185),
    @(42601.772812499999, 3078),
    @(42601.773298611108, 3258)
)

$startRow = 3
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]

    # Clone row 2 (values + styles) into the new row, then adjust the
    # two cells that differ per row.
    $ws.Range("A2:M2").Copy($ws.Range("A" + $r + ":M" + $r))
    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 3).Value = $data[1]
}

# Column A was widened to fit the refreshed content.
$ws.Range("A1").EntireColumn.ColumnWidth = 14
